$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "Consequences of ischemia most prominent"
$ws.Range("E4").Select()
